$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (M1:O1), matching the existing header style (s="1") ---
$ws.Range("L1").Copy() | Out-Null
$ws.Range("M1:O1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("M1").Value = "renewd"
$ws.Range("N1").Value = "PlanID"
$ws.Range("O1").Value = "iteration"

# --- New data columns for all data rows (2-67) ---
$ws.Range("M2:M67").Value = "after"
$ws.Range("N2:N67").Value = 20120894
$ws.Range("O2:O67").Value = 6
